$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$srcCol = $ws.Range("B2:B157")
$dstCol = $ws.Range("D2:D157")

# Move the values that live in column B over to column D (no style change
# needed -- those cells never carried one), then wipe column B so it goes
# back to being completely blank.
$dstCol.Value2 = $srcCol.Value2
$srcCol.ClearContents()

# Column D inherits the "bestFit" width that used to belong to column B;
# column B goes back to the sheet default width.
$ws.Columns("D").ColumnWidth = $ws.Columns("B").ColumnWidth
$ws.Columns("B").ColumnWidth = $ws.StandardWidth

# Matches the author's resulting selection: the entire (now data-bearing)
# column D.
$ws.Range("D1:D1048576").Select()
